# Refresh cryptocurrency prices and 1h volume-change percentages scraped from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.341.52"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "1.852.54"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "'322.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.4482"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.94%  "
$ws.Range("D8").Value = "'0.3828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.87%  "
$ws.Range("D9").Value = "'48.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.06%  "
$ws.Range("D10").Value = "'0.07864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.40%  "
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("D13").Value = "1.810.14"
$ws.Range("E13").Value = "  -6.96%  "
$ws.Range("D14").Value = "'5.860"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.50%  "
$ws.Range("D15").Value = "'7.122"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'85.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("E18").Value = "  -3.48%  "
$ws.Range("D19").Value = "'0.06506"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "'16.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.03%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "'5.466"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.89%  "
$ws.Range("D23").Value = "27.339.22"
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "'10.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.93%  "
$ws.Range("D25").Value = "'2.259"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "2.084.78"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "'151.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").Value = "'2.061"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "'5.529"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("D31").Value = "'119.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").Value = "'0.09320"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.472"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.9353"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "'3.606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").Value = "'5.248"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.06%  "
$ws.Range("D37").Value = "'0.02221"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").Value = "'0.05960"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "'8.278"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("D43").Value = "'0.1848"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'10.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.56%  "
$ws.Range("D45").Value = "'1.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.38%  "
$ws.Range("D46").Value = "'0.5637"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.09%  "
$ws.Range("D47").Value = "'12.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.64%  "
$ws.Range("D48").Value = "'1.926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.14%  "
$ws.Range("D49").Value = "'3.360"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'0.06862"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'108.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
